$d = $word.ActiveDocument

$d.Content.Find.Execute("91-79=", $true, $true, $false, $false, $false, $true, 1, $false, "38+7=", 2) | Out-Null
$d.Content.Find.Execute("51-11=", $true, $true, $false, $false, $false, $true, 1, $false, "56-1=", 2) | Out-Null
$d.Content.Find.Execute("68+8=", $true, $true, $false, $false, $false, $true, 1, $false, "0+77=", 2) | Out-Null
$d.Content.Find.Execute("68-59=", $true, $true, $false, $false, $false, $true, 1, $false, "48-18=", 2) | Out-Null
$d.Content.Find.Execute("85+5=", $true, $true, $false, $false, $false, $true, 1, $false, "98-82=", 2) | Out-Null
$d.Content.Find.Execute("28+20=", $true, $true, $false, $false, $false, $true, 1, $false, "48-10=", 2) | Out-Null
$d.Content.Find.Execute("18+59=", $true, $true, $false, $false, $false, $true, 1, $false, "42+28=", 2) | Out-Null
$d.Content.Find.Execute("51+9=", $true, $true, $false, $false, $false, $true, 1, $false, "52+40=", 2) | Out-Null
$d.Content.Find.Execute("75-22=", $true, $true, $false, $false, $false, $true, 1, $false, "41-31=", 2) | Out-Null
$d.Content.Find.Execute("21-5=", $true, $true, $false, $false, $false, $true, 1, $false, "68+9=", 2) | Out-Null
$d.Content.Find.Execute("94-67=", $true, $true, $false, $false, $false, $true, 1, $false, "19+73=", 2) | Out-Null
$d.Content.Find.Execute("40-5=", $true, $true, $false, $false, $false, $true, 1, $false, "38-16=", 2) | Out-Null
$d.Content.Find.Execute("99-49=", $true, $true, $false, $false, $false, $true, 1, $false, "35+31=", 2) | Out-Null
$d.Content.Find.Execute("29+43=", $true, $true, $false, $false, $false, $true, 1, $false, "11+12=", 2) | Out-Null
$d.Content.Find.Execute("85-36=", $true, $true, $false, $false, $false, $true, 1, $false, "26+61=", 2) | Out-Null
$d.Content.Find.Execute("50-47=", $true, $true, $false, $false, $false, $true, 1, $false, "18+49=", 2) | Out-Null
$d.Content.Find.Execute("81-44=", $true, $true, $false, $false, $false, $true, 1, $false, "20+30=", 2) | Out-Null
$d.Content.Find.Execute("35+48=", $true, $true, $false, $false, $false, $true, 1, $false, "20+17=", 2) | Out-Null
$d.Content.Find.Execute("30+60=", $true, $true, $false, $false, $false, $true, 1, $false, "91-82=", 2) | Out-Null
$d.Content.Find.Execute("58-19=", $true, $true, $false, $false, $false, $true, 1, $false, "89-45=", 2) | Out-Null
$d.Content.Find.Execute("91-73=", $true, $true, $false, $false, $false, $true, 1, $false, "21+67=", 2) | Out-Null
$d.Content.Find.Execute("32-18=", $true, $true, $false, $false, $false, $true, 1, $false, "11+50=", 2) | Out-Null
$d.Content.Find.Execute("82-26=", $true, $true, $false, $false, $false, $true, 1, $false, "88-81=", 2) | Out-Null
$d.Content.Find.Execute("49+8=", $true, $true, $false, $false, $false, $true, 1, $false, "87-2=", 2) | Out-Null
$d.Content.Find.Execute("46+34=", $true, $true, $false, $false, $false, $true, 1, $false, "34-30=", 2) | Out-Null
$d.Content.Find.Execute("71+28=", $true, $true, $false, $false, $false, $true, 1, $false, "59-43=", 2) | Out-Null
$d.Content.Find.Execute("7+91=", $true, $true, $false, $false, $false, $true, 1, $false, "69-4=", 2) | Out-Null
$d.Content.Find.Execute("17+43=", $true, $true, $false, $false, $false, $true, 1, $false, "42+53=", 2) | Out-Null
$d.Content.Find.Execute("99-55=", $true, $true, $false, $false, $false, $true, 1, $false, "73+4=", 2) | Out-Null
$d.Content.Find.Execute("32+42=", $true, $true, $false, $false, $false, $true, 1, $false, "93-46=", 2) | Out-Null
$d.Content.Find.Execute("23+55=", $true, $true, $false, $false, $false, $true, 1, $false, "88-40=", 2) | Out-Null
$d.Content.Find.Execute("13+48=", $true, $true, $false, $false, $false, $true, 1, $false, "10+26=", 2) | Out-Null
$d.Content.Find.Execute("93-2=", $true, $true, $false, $false, $false, $true, 1, $false, "73-44=", 2) | Out-Null
$d.Content.Find.Execute("31-21=", $true, $true, $false, $false, $false, $true, 1, $false, "95-16=", 2) | Out-Null
$d.Content.Find.Execute("43+15=", $true, $true, $false, $false, $false, $true, 1, $false, "43+43=", 2) | Out-Null
$d.Content.Find.Execute("44+52=", $true, $true, $false, $false, $false, $true, 1, $false, "32+44=", 2) | Out-Null
$d.Content.Find.Execute("97-4=", $true, $true, $false, $false, $false, $true, 1, $false, "34+28=", 2) | Out-Null
$d.Content.Find.Execute("29+68=", $true, $true, $false, $false, $false, $true, 1, $false, "33-30=", 2) | Out-Null
$d.Content.Find.Execute("30-24=", $true, $true, $false, $false, $false, $true, 1, $false, "87-71=", 2) | Out-Null
$d.Content.Find.Execute("21+54=", $true, $true, $false, $false, $false, $true, 1, $false, "18+40=", 2) | Out-Null
$d.Content.Find.Execute("31+20=", $true, $true, $false, $false, $false, $true, 1, $false, "66-7=", 2) | Out-Null
$d.Content.Find.Execute("58+17=", $true, $true, $false, $false, $false, $true, 1, $false, "78-19=", 2) | Out-Null
$d.Content.Find.Execute("71-69=", $true, $true, $false, $false, $false, $true, 1, $false, "0+25=", 2) | Out-Null
$d.Content.Find.Execute("59+8=", $true, $true, $false, $false, $false, $true, 1, $false, "92-23=", 2) | Out-Null
$d.Content.Find.Execute("38-0=", $true, $true, $false, $false, $false, $true, 1, $false, "4+40=", 2) | Out-Null
$d.Content.Find.Execute("14+40=", $true, $true, $false, $false, $false, $true, 1, $false, "20+61=", 2) | Out-Null
$d.Content.Find.Execute("6+72=", $true, $true, $false, $false, $false, $true, 1, $false, "46+5=", 2) | Out-Null
$d.Content.Find.Execute("73-56=", $true, $true, $false, $false, $false, $true, 1, $false, "3+61=", 2) | Out-Null
$d.Content.Find.Execute("96-57=", $true, $true, $false, $false, $false, $true, 1, $false, "13+76=", 2) | Out-Null
$d.Content.Find.Execute("37-2=", $true, $true, $false, $false, $false, $true, 1, $false, "91-67=", 2) | Out-Null
$d.Content.Find.Execute("58+22=", $true, $true, $false, $false, $false, $true, 1, $false, "1+52=", 2) | Out-Null
$d.Content.Find.Execute("66-38=", $true, $true, $false, $false, $false, $true, 1, $false, "86-57=", 2) | Out-Null
$d.Content.Find.Execute("32-23=", $true, $true, $false, $false, $false, $true, 1, $false, "23+30=", 2) | Out-Null
$d.Content.Find.Execute("28+54=", $true, $true, $false, $false, $false, $true, 1, $false, "69-38=", 2) | Out-Null
$d.Content.Find.Execute("47-12=", $true, $true, $false, $false, $false, $true, 1, $false, "24+72=", 2) | Out-Null
$d.Content.Find.Execute("26+9=", $true, $true, $false, $false, $false, $true, 1, $false, "13+38=", 2) | Out-Null
$d.Content.Find.Execute("15+67=", $true, $true, $false, $false, $false, $true, 1, $false, "42+14=", 2) | Out-Null
$d.Content.Find.Execute("94-6=", $true, $true, $false, $false, $false, $true, 1, $false, "54-46=", 2) | Out-Null
$d.Content.Find.Execute("47+0=", $true, $true, $false, $false, $false, $true, 1, $false, "33+49=", 2) | Out-Null
$d.Content.Find.Execute("24-4=", $true, $true, $false, $false, $false, $true, 1, $false, "97-93=", 2) | Out-Null
$d.Content.Find.Execute("87-78=", $true, $true, $false, $false, $false, $true, 1, $false, "3+52=", 2) | Out-Null
$d.Content.Find.Execute("41+56=", $true, $true, $false, $false, $false, $true, 1, $false, "4+58=", 2) | Out-Null
$d.Content.Find.Execute("44+16=", $true, $true, $false, $false, $false, $true, 1, $false, "96-13=", 2) | Out-Null
$d.Content.Find.Execute("53+31=", $true, $true, $false, $false, $false, $true, 1, $false, "8+89=", 2) | Out-Null
$d.Content.Find.Execute("90-8=", $true, $true, $false, $false, $false, $true, 1, $false, "38+14=", 2) | Out-Null
$d.Content.Find.Execute("87-14=", $true, $true, $false, $false, $false, $true, 1, $false, "94-39=", 2) | Out-Null
$d.Content.Find.Execute("59-18=", $true, $true, $false, $false, $false, $true, 1, $false, "93+2=", 2) | Out-Null
$d.Content.Find.Execute("69+10=", $true, $true, $false, $false, $false, $true, 1, $false, "49-0=", 2) | Out-Null
$d.Content.Find.Execute("44-11=", $true, $true, $false, $false, $false, $true, 1, $false, "82-48=", 2) | Out-Null
$d.Content.Find.Execute("78-52=", $true, $true, $false, $false, $false, $true, 1, $false, "2+45=", 2) | Out-Null
$d.Content.Find.Execute("82-19=", $true, $true, $false, $false, $false, $true, 1, $false, "91-2=", 2) | Out-Null
$d.Content.Find.Execute("48+30=", $true, $true, $false, $false, $false, $true, 1, $false, "67-51=", 2) | Out-Null
$d.Content.Find.Execute("62+20=", $true, $true, $false, $false, $false, $true, 1, $false, "34+53=", 2) | Out-Null
$d.Content.Find.Execute("91-88=", $true, $true, $false, $false, $false, $true, 1, $false, "23-6=", 2) | Out-Null
$d.Content.Find.Execute("83-49=", $true, $true, $false, $false, $false, $true, 1, $false, "11+86=", 2) | Out-Null
$d.Content.Find.Execute("98-64=", $true, $true, $false, $false, $false, $true, 1, $false, "52-38=", 2) | Out-Null
$d.Content.Find.Execute("28+48=", $true, $true, $false, $false, $false, $true, 1, $false, "66+25=", 2) | Out-Null
$d.Content.Find.Execute("72-18=", $true, $true, $false, $false, $false, $true, 1, $false, "91-40=", 2) | Out-Null
$d.Content.Find.Execute("69+15=", $true, $true, $false, $false, $false, $true, 1, $false, "76-51=", 2) | Out-Null
$d.Content.Find.Execute("12+35=", $true, $true, $false, $false, $false, $true, 1, $false, "31-6=", 2) | Out-Null
$d.Content.Find.Execute("23+39=", $true, $true, $false, $false, $false, $true, 1, $false, "77-47=", 2) | Out-Null
$d.Content.Find.Execute("68-55=", $true, $true, $false, $false, $false, $true, 1, $false, "62-9=", 2) | Out-Null
$d.Content.Find.Execute("81-21=", $true, $true, $false, $false, $false, $true, 1, $false, "39+46=", 2) | Out-Null
$d.Content.Find.Execute("45+19=", $true, $true, $false, $false, $false, $true, 1, $false, "34-31=", 2) | Out-Null
$d.Content.Find.Execute("53-16=", $true, $true, $false, $false, $false, $true, 1, $false, "23+13=", 2) | Out-Null
$d.Content.Find.Execute("12+37=", $true, $true, $false, $false, $false, $true, 1, $false, "22+18=", 2) | Out-Null
$d.Content.Find.Execute("21+71=", $true, $true, $false, $false, $false, $true, 1, $false, "62-50=", 2) | Out-Null
$d.Content.Find.Execute("39+28=", $true, $true, $false, $false, $false, $true, 1, $false, "54-20=", 2) | Out-Null
$d.Content.Find.Execute("82+8=", $true, $true, $false, $false, $false, $true, 1, $false, "1+89=", 2) | Out-Null
$d.Content.Find.Execute("24+29=", $true, $true, $false, $false, $false, $true, 1, $false, "1+33=", 2) | Out-Null
$d.Content.Find.Execute("91-5=", $true, $true, $false, $false, $false, $true, 1, $false, "57+34=", 2) | Out-Null
$d.Content.Find.Execute("83+16=", $true, $true, $false, $false, $false, $true, 1, $false, "54-4=", 2) | Out-Null
$d.Content.Find.Execute("69-51=", $true, $true, $false, $false, $false, $true, 1, $false, "42-18=", 2) | Out-Null
$d.Content.Find.Execute("9+18=", $true, $true, $false, $false, $false, $true, 1, $false, "76+16=", 2) | Out-Null
$d.Content.Find.Execute("10+66=", $true, $true, $false, $false, $false, $true, 1, $false, "78-62=", 2) | Out-Null
$d.Content.Find.Execute("73-42=", $true, $true, $false, $false, $false, $true, 1, $false, "14+34=", 2) | Out-Null
$d.Content.Find.Execute("33+5=", $true, $true, $false, $false, $false, $true, 1, $false, "63-28=", 2) | Out-Null
$d.Content.Find.Execute("2+26=", $true, $true, $false, $false, $false, $true, 1, $false, "82-31=", 2) | Out-Null
$d.Content.Find.Execute("63-11=", $true, $true, $false, $false, $false, $true, 1, $false, "12+84=", 2) | Out-Null
$d.Content.Find.Execute("93-58=", $true, $true, $false, $false, $false, $true, 1, $false, "93-5=", 2) | Out-Null
